$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rewrite the label/value cells that already existed in rows 10-21 ---
$ws.Range("A10").Value = "Objetivos:"
$ws.Range("B10").Value = "8767640 - Eduardo Ferro dos Santos"
$ws.Range("C10").Value = "8767640 - Eduardo Ferro dos Santos"
$ws.Range("A11").Value = "Objectives:"
$ws.Range("A12").Value = "Docentes responsáveis:"
$ws.Range("A15").Value = "Programa:"
$ws.Range("B15").Value = "8767640 - Eduardo Ferro dos Santos"
$ws.Range("C15").Value = "8767640 - Eduardo Ferro dos Santos"
$ws.Range("A16").Value = "Syllabus:"
$ws.Range("A17").Value = "Avaliação:"
$ws.Range("A18").Value = "Método:"
$ws.Range("A19").Value = "Critério:"
$ws.Range("A20").Value = "Norma de recuperação:"
$ws.Range("B20").Value = "Nota Final = Média Ponderada das formas de avaliação definidas pelo docente, abrangendo avaliações individuais e em grupo."
$ws.Range("C20").Value = "Nota Final = Média Ponderada das formas de avaliação definidas pelo docente, abrangendo avaliações individuais e em grupo."
$ws.Range("A21").Value = "Bibliografia:"
$ws.Range("B21").Value = "NF = (MF + PR)/2, onde NF é a média final da segunda avaliação, MF é a média final da primeira avaliação e PR é a nota da recuperação."
$ws.Range("C21").Value = "NF = (MF + PR)/2, onde NF é a média final da segunda avaliação, MF é a média final da primeira avaliação e PR é a nota da recuperação."

# --- B13/C13 become the text "01/01/2021"; assign it as a values-only paste from
#     B8/C8 (same literal text) so Excel keeps it as text instead of parsing a date ---
$ws.Range("B8").Copy()
$ws.Range("B13").PasteSpecial(-4163)
$ws.Range("C8").Copy()
$ws.Range("C13").PasteSpecial(-4163)

# --- New cells: copy formatting from a same-column neighbor first, then set the value ---
$ws.Range("A15").Copy()
$ws.Range("A13").PasteSpecial(-4122)
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("A15").Copy()
$ws.Range("A14").PasteSpecial(-4122)
$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("B15").Copy()
$ws.Range("B18").PasteSpecial(-4122)
$ws.Range("B18").Value = "8188658 - Maria Auxiliadora Motta Barreto"
$ws.Range("C15").Copy()
$ws.Range("C18").PasteSpecial(-4122)
$ws.Range("C18").Value = "8188658 - Maria Auxiliadora Motta Barreto"
$ws.Range("B15").Copy()
$ws.Range("B19").PasteSpecial(-4122)
$ws.Range("B19").Value = "Aulas Expositivas, Aulas Baseadas em Problemas e Projetos, Atividades Individuais e em Grupo, Seminários. Dadas estas características, haverá múltiplas formas de avaliação definidas pelo docente."
$ws.Range("C15").Copy()
$ws.Range("C19").PasteSpecial(-4122)
$ws.Range("C19").Value = "Aulas Expositivas, Aulas Baseadas em Problemas e Projetos, Atividades Individuais e em Grupo, Seminários. Dadas estas características, haverá múltiplas formas de avaliação definidas pelo docente."

# --- Cells that no longer exist in the rebuilt rows are fully removed ---
$ws.Range("B14").Clear()
$ws.Range("C14").Clear()
$ws.Range("B17").Clear()
$ws.Range("C17").Clear()

# --- Apply the custom row heights required by the new layout ---
$ws.Rows.Item(13).RowHeight = 60
$ws.Rows.Item(14).RowHeight = 60
$ws.Rows.Item(15).RowHeight = 120
$ws.Rows.Item(16).RowHeight = 120
$ws.Rows.Item(18).RowHeight = 60
$ws.Rows.Item(19).RowHeight = 60
$ws.Rows.Item(21).RowHeight = 120

# --- Row 17 returns to the sheet default height (no more custom height) ---
$ws.Rows.Item(17).EntireRow.AutoFit()

# --- Remove the two trailing rows that no longer exist in the rebuilt table ---
$ws.Rows.Item(23).EntireRow.Delete()
$ws.Rows.Item(22).EntireRow.Delete()

$excel.CutCopyMode = 0
